$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing rows 9 and 10: B column value changes from "Informatique" to "SV" ---
$ws.Range("B9").Value = "SV"
$ws.Range("B10").Value = "SV"

# --- Add new rows 12-19, written column-by-column so new shared strings are
#     interned in the same order as the source workbook (B, then C, then D) ---
$colA = @(2025, 2025, 2025, 2025, 2025, 2025, 2025, 2025)
$colB = @("SV", "CH", "CH", "PUYSQ", "CH", "PUYSQ", "CH", "PUYSQ")
$colC = @("L2", "L3", "L2", "L2", "L3", "L2", "L2", "L3")
$colD = @("LIM1", "LIM2", "LIM3", "LIM4", "LIM5", "LIM6", "LIM7", "LIM8")
$colE = @("M", "F", "M", "F", "M", "F", "M", "F")
$colF = @(200, 180, 150, 22, 32, 45, 45, 40)

for ($i = 0; $i -lt 8; $i++) {
    $ws.Cells.Item(12 + $i, 1).Value = $colA[$i]
}
for ($i = 0; $i -lt 8; $i++) {
    $ws.Cells.Item(12 + $i, 2).Value = $colB[$i]
}
for ($i = 0; $i -lt 8; $i++) {
    $ws.Cells.Item(12 + $i, 3).Value = $colC[$i]
}
for ($i = 0; $i -lt 8; $i++) {
    $ws.Cells.Item(12 + $i, 4).Value = $colD[$i]
}
for ($i = 0; $i -lt 8; $i++) {
    $ws.Cells.Item(12 + $i, 5).Value = $colE[$i]
}
for ($i = 0; $i -lt 8; $i++) {
    $ws.Cells.Item(12 + $i, 6).Value = $colF[$i]
}

# Apply the same formatting as the rest of the data rows (vertical-center + wrap text)
$ws.Range("A2").Copy()
$ws.Range("A12:F19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Register the small (8pt) Times New Roman font used for the phonetic-guide
#     default, matching the font table of the target workbook ---
$tmp = $ws.Range("H1")
$tmp.Value = "x"
$tmp.Font.Size = 8
$tmp.Clear()

# --- Selection moves to J20, as in the target file ---
[void]$ws.Range("J20").Select()
